$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": M2 (PORCELANATO) 133.36 -> 1410.75
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 1410.75

# Sheet "VENTA MENSUAL": F2 (junio) 133.36 -> 1410.75, F30 total 2524.59 -> 3801.98
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 1410.75
$ws2.Range("F30").Value = 3801.98
$ws2.Columns.Item(6).ColumnWidth = 12.166666666666666

# Sheet "CUMPLIMIENTO MENSUAL": row 16 (PORCELANATO) and row 19 (TOTAL)
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 2757.31
$ws3.Range("E16").Value = 16041.3
$ws3.Range("F16").Value = 0.1466762702135956
$ws3.Range("D19").Value = 3796.22
$ws3.Range("E19").Value = 25741.57107555787
$ws3.Range("F19").Value = 0.1285207817432672
$ws3.Columns.Item(6).ColumnWidth = 23.166666666666668
